# Fix the CLAVE (course code) column (A) so every row of a course's block
# shows the correct course code instead of the stray per-topic placeholder
# strings that used to live there. Also corrects the mid-block code change
# for "Impuestos Sobre la Renta personas morales" (rows 64-66 really belong
# to CFS27606, not CFS27603) and normalizes "CAU 27401" -> "CAU27401".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$TAB = [char]9

$groups = @(
    @{ Start = 3;  End = 8;  Code = "CAU27401" + $TAB },
    @{ Start = 9;  End = 14; Code = "CCE27401" },
    @{ Start = 15; End = 20; Code = "CCN27401 " },
    @{ Start = 21; End = 26; Code = "CAF27401" },
    @{ Start = 27; End = 32; Code = "CFS27401" },
    @{ Start = 61; End = 63; Code = "CFS27603" },
    @{ Start = 64; End = 66; Code = "CFS27606" },
    @{ Start = 67; End = 72; Code = "CAF27603" },
    @{ Start = 73; End = 78; Code = "CCN27603        " },
    @{ Start = 79; End = 84; Code = "CCE27603" }
)

foreach ($g in $groups) {
    for ($r = $g.Start; $r -le $g.End; $r++) {
        $ws.Range("A$r").Value = $g.Code
    }
}

$ws.Range("A79:A84").Select()
